$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.037.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "'2.285.31"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.20%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'316.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").Value = "'102.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.33%  "

$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.603"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("E10").Value = "  -4.39%  "

$ws.Range("D11").Value = "'0.0905"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").Value = "'8.30"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("D13").Value = "'0.106"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").Value = "'0.964"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.33%  "

$ws.Range("D15").Value = "'15.26"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.76%  "

$ws.Range("D16").Value = "'2.633.14"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.28%  "

$ws.Range("D17").Value = "'2.287.51"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.91%  "

$ws.Range("D18").Value = "'41.940.91"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").Value = "'7.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").Value = "'283.63"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +10.96%  "

$ws.Range("D22").Value = "'73.62"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("D23").Value = "'3.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.86%  "

$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").Value = "'9.87"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.79%  "

$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").Value = "'10.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.83%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.27"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.08%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'23.11"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.69%  "

$ws.Range("D30").Value = "'163.38"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.94%  "

$ws.Range("D31").Value = "'34.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.53%  "

$ws.Range("D32").Value = "'0.0876"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("D34").Value = "'5.83"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.51%  "

$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("E36").Value = "  -7.49%  "

$ws.Range("D37").Value = "'4.57"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("E38").Value = "  +9.58%  "

$ws.Range("E39").Value = "  -3.45%  "

$ws.Range("D40").Value = "'3.59"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.18%  "

$ws.Range("D41").Value = "'102.33"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +20.63%  "

$ws.Range("D42").Value = "'1.46"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").Value = "'69.66"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.87%  "

$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("E45").Value = "  -4.04%  "

$ws.Range("D46").Value = "'115.58"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.29%  "

$ws.Range("D47").Value = "'11.94"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("D48").Value = "'9.03"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("D49").Value = "'76.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.99%  "

$ws.Range("E50").Value = "  -2.02%  "

$ws.Range("D51").Value = "'1.26"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.68%  "
